# Trade #113 closed at 2026-02-16 21:43:17 - leadlag UP +0.000%
#
# This commit:
#  1. Opens a new leadlag trade (#113) -> new row in "leadlag" sheet + "All Trades" summary row not needed
#     (All Trades sheet only logs closed-state updates here; new OPEN leadlag trade only appears on its
#     own strategy sheet, matching source diff).
#  2. Closes two previously-OPEN momentum trades (#89 and #90) -> updates rows 24/25 on "momentum" sheet,
#     and appends their now-CLOSED state as new rows on "All Trades".
#  3. Refreshes the roll-up stats on "Summary" and "Comparison" sheets to reflect the above.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary": update OVERALL and momentum rows with new trade counts.
# ---------------------------------------------------------------------------
$sum = $wb.Worksheets.Item("Summary")

$sum.Range("C2").Value = 90
$sum.Range("D2").Value = "'71.1%"
$sum.Range("E2").Value = "'+26.8326%"
$sum.Range("F2").Value = "'+0.2981%"

$sum.Range("C4").Value = 26
$sum.Range("D4").Value = "'84.6%"
$sum.Range("E4").Value = "'+12.6854%"
$sum.Range("F4").Value = "'+0.4879%"

# ---------------------------------------------------------------------------
# Sheet "leadlag": append newly opened trade #113.
# ---------------------------------------------------------------------------
$lead = $wb.Worksheets.Item("leadlag")

$lead.Cells.Item(88, 1).Value = 113
$lead.Cells.Item(88, 2).Value = "'2026-02-16"
$lead.Cells.Item(88, 3).Value = "21:43:17"
$lead.Cells.Item(88, 4).Value = "leadlag"
$lead.Cells.Item(88, 5).Value = "UP"
$lead.Cells.Item(88, 6).Value = 68404.295
$lead.Cells.Item(88, 8).Value = "OPEN"
$lead.Cells.Item(88, 9).Value = 0
$lead.Cells.Item(88, 10).Value = 0
$lead.Cells.Item(88, 11).Value = 0.75
$lead.Cells.Item(88, 12).Value = "Coinbase leading with 0.078% move"
$lead.Cells.Item(88, 14).Value = 0

# ---------------------------------------------------------------------------
# Sheet "momentum": close out trades #89 (row 24) and #90 (row 25).
# ---------------------------------------------------------------------------
$mom = $wb.Worksheets.Item("momentum")

$mom.Cells.Item(24, 7).Value = 67918.949609
$mom.Cells.Item(24, 8).Value = "CLOSED"
$mom.Cells.Item(24, 9).Value = 0.7665
$mom.Cells.Item(24, 10).Value = 7.67
$mom.Cells.Item(24, 13).Value = "time_exit_5min"
$mom.Cells.Item(24, 14).Value = 5

$mom.Cells.Item(25, 7).Value = 68303.717177
$mom.Cells.Item(25, 8).Value = "CLOSED"
$mom.Cells.Item(25, 9).Value = 0.2353
$mom.Cells.Item(25, 10).Value = 2.35
$mom.Cells.Item(25, 13).Value = "time_exit_5min"
$mom.Cells.Item(25, 14).Value = 5

# ---------------------------------------------------------------------------
# Sheet "All Trades": append the two momentum trades in their closed state.
# ---------------------------------------------------------------------------
$all = $wb.Worksheets.Item("All Trades")

$all.Cells.Item(90, 1).Value = 89
$all.Cells.Item(90, 2).Value = "'2026-02-16"
$all.Cells.Item(90, 3).Value = "21:38:07"
$all.Cells.Item(90, 4).Value = "momentum"
$all.Cells.Item(90, 5).Value = "DOWN"
$all.Cells.Item(90, 6).Value = 68443.59
$all.Cells.Item(90, 7).Value = 67918.949609
$all.Cells.Item(90, 8).Value = "CLOSED"
$all.Cells.Item(90, 9).Value = 0.7665
$all.Cells.Item(90, 10).Value = 7.67
$all.Cells.Item(90, 11).Value = 0.9
$all.Cells.Item(90, 12).Value = "Downward momentum: -0.225% over 10 samples"
$all.Cells.Item(90, 13).Value = "time_exit_5min"
$all.Cells.Item(90, 14).Value = 5

$all.Cells.Item(91, 1).Value = 90
$all.Cells.Item(91, 2).Value = "'2026-02-16"
$all.Cells.Item(91, 3).Value = "21:38:13"
$all.Cells.Item(91, 4).Value = "momentum"
$all.Cells.Item(91, 5).Value = "DOWN"
$all.Cells.Item(91, 6).Value = 68464.83500000001
$all.Cells.Item(91, 7).Value = 68303.717177
$all.Cells.Item(91, 8).Value = "CLOSED"
$all.Cells.Item(91, 9).Value = 0.2353
$all.Cells.Item(91, 10).Value = 2.35
$all.Cells.Item(91, 11).Value = 0.9
$all.Cells.Item(91, 12).Value = "Downward momentum: -0.243% over 10 samples"
$all.Cells.Item(91, 13).Value = "time_exit_5min"
$all.Cells.Item(91, 14).Value = 5

# ---------------------------------------------------------------------------
# Sheet "Comparison": refresh momentum strategy rollup row.
# ---------------------------------------------------------------------------
$comp = $wb.Worksheets.Item("Comparison")

$comp.Range("B3").Value = 26
$comp.Range("C3").Value = "'84.6%"
$comp.Range("D3").Value = "'12.28"
$comp.Range("E3").Value = "'+0.6277%"
$comp.Range("G3").Value = "'1.12"
